$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to update so that numeric-
# and percentage-looking strings (e.g. "301.73", "-0.64%") are stored as
# literal text, matching the worksheet's existing inline-string convention
# instead of being auto-converted to Number/Percentage by Excel.
$updatedCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $updatedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the refreshed price / 1h-volume values scraped on 2023-01-22.
$ws.Range("D2").Value = "301.73"
$ws.Range("E2").Value = "-0.64%"
$ws.Range("D3").Value = "37.55"
$ws.Range("E3").Value = "7.63%"
$ws.Range("D4").Value = "5.006"
$ws.Range("E4").Value = "-2.42%"
$ws.Range("D5").Value = "0.07856"
$ws.Range("E5").Value = "1.22%"
$ws.Range("D6").Value = "2.234"
$ws.Range("E6").Value = "-5.53%"
$ws.Range("D7").Value = "8.020"
$ws.Range("E7").Value = "-0.10%"
$ws.Range("D8").Value = "4.025"
$ws.Range("E8").Value = "2.32%"
$ws.Range("D9").Value = "0.9106"
$ws.Range("E9").Value = "-1.98%"
$ws.Range("D10").Value = "0.1882"
$ws.Range("E10").Value = "4.96%"
$ws.Range("D11").Value = "0.09231"
$ws.Range("E11").Value = "-9.06%"
$ws.Range("D12").Value = "0.08524"
$ws.Range("E12").Value = "-0.13%"
$ws.Range("D13").Value = "0.03531"
$ws.Range("E13").Value = "6.65%"
$ws.Range("D14").Value = "0.09932"
$ws.Range("E14").Value = "0.43%"
$ws.Range("D15").Value = "0.001486"
$ws.Range("E15").Value = "-0.65%"
$ws.Range("D16").Value = "0.005664"
$ws.Range("E16").Value = "-1.68%"
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").Value = "0.16%"
$ws.Range("E18").Value = "0.99%"
$ws.Range("E19").Value = "2.85%"
$ws.Range("D20").Value = "0.1307"
$ws.Range("E20").Value = "-0.26%"
$ws.Range("D21").Value = "4.779"
$ws.Range("E21").Value = "10.94%"
$ws.Range("D22").Value = "0.2200"
$ws.Range("E23").Value = "1.79%"
$ws.Range("E24").Value = "0.90%"
$ws.Range("D25").Value = "0.004450"
$ws.Range("E25").Value = "-0.23%"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").Value = "-0.07%"
$ws.Range("D27").Value = "0.0004741"
$ws.Range("E27").Value = "28.24%"
$ws.Range("D39").Value = "0.01769"
$ws.Range("E39").Value = "-1.09%"
$ws.Range("D40").Value = "0.04725"
$ws.Range("E40").Value = "-0.60%"
$ws.Range("D41").Value = "0.007834"
$ws.Range("E41").Value = "1.27%"
$ws.Range("D42").Value = "0.1392"
$ws.Range("E42").Value = "-1.39%"
$ws.Range("E43").Value = "7.99%"
$ws.Range("D44").Value = "0.002216"
$ws.Range("E44").Value = "5.32%"
$ws.Range("D45").Value = "0.01024"
$ws.Range("E45").Value = "11.30%"
$ws.Range("D46").Value = "0.00005982"
$ws.Range("E46").Value = "-2.18%"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("D48").Value = "8.671"
$ws.Range("E48").Value = "218.12%"
$ws.Range("D49").Value = "0.002685"
$ws.Range("E49").Value = "34.34%"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.09%"
